# Auto-generated script applying the cryptos.xlsx price/volume update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.806.44"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.87%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.559.46"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.55"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "188.40"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.94%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.553.45"
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.84%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.620"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.60%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.213"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +6.48%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.645"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.00"
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000307"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.42"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.136.59"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "70.813.21"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.93%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.572.43"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.95%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.98"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.40%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.66"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.10%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "568.45"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.51%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.120"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.995"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.43%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.75"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.79%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.58"
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.89"
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "93.84"
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.84%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.10"
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.91"
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.45%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.27"
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.79%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.38"
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.25%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.06"
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.09%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.19"
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.10%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("B33").Style = "Normal"
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "dogwifhat"
$ws.Range("B33").Style = "Normal"

$ws.Range("C33").Style = "Normal"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C33").Style = "Normal"

$ws.Range("D33").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.92"
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +23.71%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("B34").Style = "Normal"
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Hedera"
$ws.Range("B34").Style = "Normal"

$ws.Range("C34").Style = "Normal"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C34").Style = "Normal"

$ws.Range("D34").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.115"
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "63.09"
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.34%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.24"
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.61%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("B37").Style = "Normal"
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "TheGraph"
$ws.Range("B37").Style = "Normal"

$ws.Range("C37").Style = "Normal"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("C37").Style = "Normal"

$ws.Range("D37").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.407"
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").Style = "Normal"
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "Bittensor"
$ws.Range("B38").Style = "Normal"

$ws.Range("C38").Style = "Normal"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C38").Style = "Normal"

$ws.Range("D38").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "525.91"
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.16%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.19"
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.30%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.622.48"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +8.54%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0786"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.72%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("B43").Style = "Normal"
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Stacks"
$ws.Range("B43").Style = "Normal"

$ws.Range("C43").Style = "Normal"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C43").Style = "Normal"

$ws.Range("D43").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.53"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.34%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("B44").Style = "Normal"
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Kaspa"
$ws.Range("B44").Style = "Normal"

$ws.Range("C44").Style = "Normal"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C44").Style = "Normal"

$ws.Range("D44").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.139"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.79%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0458"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.09%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Style = "Normal"
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("B46").Style = "Normal"

$ws.Range("C46").Style = "Normal"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("C46").Style = "Normal"

$ws.Range("D46").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.45"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.33%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Style = "Normal"
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("B47").Style = "Normal"

$ws.Range("C47").Style = "Normal"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("C47").Style = "Normal"

$ws.Range("D47").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.92"
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.77%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.137"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.64%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.19"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.92%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +8.98%  "
$ws.Range("E51").Style = "Normal"

